$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Add a new row to the "Language" table for Lisp.
$newRow = $tbl.ListRows.Add()
$newRow.Range.Cells.Item(1, 1).Value = "Lisp"
$newRow.Range.Cells.Item(1, 2).Value = "Yes"
$newRow.Range.Cells.Item(1, 3).Value = "Yes"
$newRow.Range.Cells.Item(1, 4).Value = 1958
$newRow.Range.Cells.Item(1, 5).Value = "VM"
$newRow.Range.Cells.Item(1, 6).Value = 15199
$newRow.Range.Cells.Item(1, 7).Value = "Multi"
$newRow.Range.Cells.Item(1, 8).Value = "Global"

# The table is kept sorted by GitHub Repos (column F) ascending; re-sort
# now that the new row has been appended so it lands in the right place.
$tbl.Sort.SortFields.Clear()
$tbl.Sort.SortFields.Add($tbl.ListColumns.Item("GitHub Repos").Range)
$tbl.Sort.Header = 1
$tbl.Sort.Apply()

# Leave the selection where the user clicked after adding the row.
[void]$ws.Range("A18").Select()
